# Update sfbic_model.xlsx with refreshed regression output
# (re-fit logistic model: updated predictor list + new rows for the
#  two additional predictors that became significant).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Estatisticas" - refreshed model fit statistics (row 2)
# ---------------------------------------------------------------------
$wsStats = $wb.Worksheets.Item("Estatisticas")

$wsStats.Range("A2").Value = -137.238027615703
$wsStats.Range("B2").Value = 296.476055231405
$wsStats.Range("C2").Value = 344.915371408688
$wsStats.Range("E2").Value = 0.91887417218543
$wsStats.Range("F2").Value = 0.982824427480916
$wsStats.Range("G2").Value = 0.5
$wsStats.Range("H2").Value = 0.927927927927928
$wsStats.Range("I2").Value = 0.816326530612245

# ---------------------------------------------------------------------
# Shared predictor labels (column A) used by Coeficientes / Significativas
# ---------------------------------------------------------------------
$labels = @(
    "(Intercept)",
    "``De forma geral, como você avalia o nível de suporte recebido do seu orientador durante a realização de sua dissertação/tese de pós-graduação?``",
    "``Em que ano você iniciou seu último curso de pós-graduação stricto sensu na UFMT?``",
    "``Especifique a natureza da jornada de trabalho durante a pós-graduação_Jornada de tempo integral (acima de 30 horas semanais)``",
    "``Como você classificaria seu nível de desempenho acadêmico ao longo do curso de pós-graduação?``",
    "``Você fez uso de algum período de prorrogação de prazo para o término do seu curso de pós-graduação? Se sim, qual foi o período de prorrogação utilizado?_Não prorroguei o período de curso``",
    "``Qual a probabilidade de você recomendar esta instituição para futuros estudantes de pós-graduação?``",
    "``Entre os fatores abaixo, quais você acredita que poderiam ser aprimorados para melhorar significativamente a sua experiência com a instituição? Escolha todos que se aplicam_Integração da instituição com a comunidade local``",
    "``Você trancou sua matrícula no programa de pós-graduação? Se sim, por quanto tempo?_Não tranquei a matrícula``",
    "``Você trabalhou enquanto estava cursando pós-graduação?_Sim, durante todo o período do curso``",
    "``Qual era o rendimento per capita familiar durante o período em que realizou o curso de pós-graduação?_Sem renda``"
)

# Coefficient table: Estimate, Std. Error, z value, Pr(>|z|) for rows 2..12
$coef = @(
    @(-483.631365048737, 98.5331784571956, -4.90830979596212, 0.000000918646628247724),
    @(-0.740910385420564, 0.138871581818701, -5.33521960157287, 0.0000000954289920562034),
    @(0.241996372352355, 0.0489212881299963, 4.94664759663133, 0.000000755024837734198),
    @(1.10176160471818, 0.36894822209533, 2.98622283219324, 0.00282446747986601),
    @(-0.69416287492339, 0.190754584707591, -3.63903638797189, 0.000273660132916959),
    @(1.18625009595802, 0.364845615496839, 3.25137550123113, 0.0011484806284036),
    @(-0.399024425782989, 0.147475723620061, -2.7056956629078, 0.00681614808082955),
    @(-1.03205147766695, 0.346492756061415, -2.97856581302964, 0.00289600828795056),
    @(-1.36228918171106, 0.639050683610116, -2.13173886931039, 0.0330283195513054),
    @(1.17656683858954, 0.388032612224987, 3.03213390194985, 0.00242831458431084),
    @(2.16559694620862, 0.7951093615082, 2.72364664667111, 0.00645655443835115)
)

# Odds-ratio table: Estimate, OR, CI_lower, CI_upper for rows 2..12
# (row 2's OR/CI values are extreme-magnitude doubles; PowerShell's parser
#  here has no literal `e` notation, so build them via a [double] string cast)
$row2OR = [double]"9.15307211314962e-211"
$row2CIlo = [double]"1.90222172835705e-299"
$row2CIhi = [double]"2.06338816366649e-131"

$odds = @(
    @(-483.631365048737, $row2OR, $row2CIlo, $row2CIhi),
    @(-0.740910385420564, 0.476679755624802, 0.359849743497069, 0.621736088940748),
    @(0.241996372352355, 1.27378957194807, 1.16333873081158, 1.40972367943797),
    @(1.10176160471818, 3.00946284106764, 1.4771315389662, 6.3118769063912),
    @(-0.69416287492339, 0.499492410639741, 0.339532268659008, 0.719448809242678),
    @(1.18625009595802, 3.27477805052307, 1.64274211448891, 6.91849053657683),
    @(-0.399024425782989, 0.670974312080348, 0.500369401550011, 0.893952632236859),
    @(-1.03205147766695, 0.356275319491855, 0.176167584306328, 0.689301019377057),
    @(-1.36228918171106, 0.25607390578067, 0.0717735812157048, 0.88878763739367),
    @(1.17656683858954, 3.2432205680308, 1.53947300053348, 7.09912922572255),
    @(2.16559694620862, 8.71980561843554, 1.72554246231028, 40.6800152608079)
)

function Update-CoefficientSheet($ws) {
    for ($i = 0; $i -lt $labels.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $labels[$i]
    }
    for ($i = 0; $i -lt $coef.Count; $i++) {
        $row = $i + 2
        $vals = $coef[$i]
        $ws.Cells.Item($row, 2).Value = $vals[0]
        $ws.Cells.Item($row, 3).Value = $vals[1]
        $ws.Cells.Item($row, 4).Value = $vals[2]
        $ws.Cells.Item($row, 5).Value = $vals[3]
    }
}

# ---------------------------------------------------------------------
# Sheet "Coeficientes"
# ---------------------------------------------------------------------
$wsCoef = $wb.Worksheets.Item("Coeficientes")
Update-CoefficientSheet $wsCoef

# ---------------------------------------------------------------------
# Sheet "Significativas" (mirrors Coeficientes)
# ---------------------------------------------------------------------
$wsSig = $wb.Worksheets.Item("Significativas")
Update-CoefficientSheet $wsSig

# ---------------------------------------------------------------------
# Sheet "Odds Ratios"
# ---------------------------------------------------------------------
$wsOR = $wb.Worksheets.Item("Odds Ratios")

for ($i = 0; $i -lt $labels.Count; $i++) {
    $row = $i + 2
    $wsOR.Cells.Item($row, 1).Value = $labels[$i]
}
for ($i = 0; $i -lt $odds.Count; $i++) {
    $row = $i + 2
    $vals = $odds[$i]
    $wsOR.Cells.Item($row, 2).Value = $vals[0]
    $wsOR.Cells.Item($row, 3).Value = $vals[1]
    $wsOR.Cells.Item($row, 4).Value = $vals[2]
    $wsOR.Cells.Item($row, 5).Value = $vals[3]
}
